# (1) Rename the original sheet "RData" -> "Cost"
$wb = $excel.ActiveWorkbook
$wsCost = $wb.Worksheets.Item(1)
$wsCost.Name = "Cost"

# (2) Add a new "Lookup" sheet right after "Cost"
$wsLookup = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCost)
$wsLookup.Name = "Lookup"

# Fill in the lookup table, column by column (A, then B, then C) so that new
# shared strings get appended to the shared-string table in the same order
# the original authoring tool produced them in.

# --- Column A ---
$wsLookup.Range("A1").Value = "sname"
$wsLookup.Range("A2").Value = "cdmards"
$wsLookup.Range("A3").Value = "abtivmtx"
$wsLookup.Range("A4").Value = "adamtx"
$wsLookup.Range("A5").Value = "ada"
$wsLookup.Range("A6").Value = "etnmtx"
$wsLookup.Range("A7").Value = "etn"
$wsLookup.Range("A8").Value = "golmtx"
$wsLookup.Range("A9").Value = "ifxmtx"
$wsLookup.Range("A10").Value = "placebo"
$wsLookup.Range("A11").Value = "tczmtx"
$wsLookup.Range("A12").Value = "tcz"
$wsLookup.Range("A13").Value = "czpmtx"
$wsLookup.Range("A14").Value = "abtscmtx"
$wsLookup.Range("A15").Value = "nbt"
$wsLookup.Range("A16").Value = "rtxmtx"
$wsLookup.Range("A17").Value = "tofmtx"
$wsLookup.Range("A18").Value = "rtx"
$wsLookup.Range("A19").Value = "tof"
$wsLookup.Range("A20").Value = "czp"
$wsLookup.Range("A21").Value = "gol"

# --- Column B ---
$wsLookup.Range("B1").Value = "agent1"
$wsLookup.Range("B2").Value = "cdmards"
$wsLookup.Range("B3").Value = "abtiv"
$wsLookup.Range("B4").Value = "ada"
$wsLookup.Range("B5").Value = "ada"
$wsLookup.Range("B6").Value = "etn"
$wsLookup.Range("B7").Value = "etn"
$wsLookup.Range("B8").Value = "gol"
$wsLookup.Range("B9").Value = "ifx"
$wsLookup.Range("B10").Value = "placebo"
$wsLookup.Range("B11").Value = "tcz"
$wsLookup.Range("B12").Value = "tcz"
$wsLookup.Range("B13").Value = "czp"
$wsLookup.Range("B14").Value = "abtsc"
$wsLookup.Range("B15").Value = "nbt"
$wsLookup.Range("B16").Value = "rtx"
$wsLookup.Range("B17").Value = "tof"
$wsLookup.Range("B18").Value = "rtx"
$wsLookup.Range("B19").Value = "tof"
$wsLookup.Range("B20").Value = "czp"
$wsLookup.Range("B21").Value = "gol"

# --- Column C ---
$wsLookup.Range("C1").Value = "agent2"
$wsLookup.Range("C3").Value = "cdmards"
$wsLookup.Range("C4").Value = "cdmards"
$wsLookup.Range("C6").Value = "cdmards"
$wsLookup.Range("C8").Value = "cdmards"
$wsLookup.Range("C9").Value = "cdmards"
$wsLookup.Range("C11").Value = "cdmards"
$wsLookup.Range("C13").Value = "cdmards"
$wsLookup.Range("C14").Value = "cdmards"
$wsLookup.Range("C16").Value = "cdmards"
$wsLookup.Range("C17").Value = "cdmards"

# Match the page setup used on the rest of the workbook's sheets
$psLookup = $wsLookup.PageSetup
$psLookup.LeftMargin = 0.75 * 72
$psLookup.RightMargin = 0.75 * 72
$psLookup.TopMargin = 1 * 72
$psLookup.BottomMargin = 1 * 72
$psLookup.HeaderMargin = 0.5 * 72
$psLookup.FooterMargin = 0.5 * 72
$psLookup.Orientation = 1

# Selection on the Lookup sheet ends up on C20 (matches authoring tool state)
$null = $wsLookup.Range("C20").Select()

# Selection on the Cost sheet ends up on C3
$null = $wsCost.Select()
$null = $wsCost.Range("C3").Select()
